$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "total marks" calculation error: Right-answer marking went from 5 to 4,
# and Wrong-answer penalty went from -1 to -2, changing the totals.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "78 / 112"
